$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 161.3069735234315
$ws.Range("C3").Value = 11.41176323146423
$ws.Range("C4").Value = 8.936350316468626
$ws.Range("C5").Value = 14.45465633022249
$ws.Range("C6").Value = 35.50115077440965
$ws.Range("C7").Value = 11.08122065020042
$ws.Range("C8").Value = 8.616045779757231
$ws.Range("C9").Value = 26.34102605558101
$ws.Range("C10").Value = 41.28564594652642
$ws.Range("C11").Value = 8.673086313692135
$ws.Range("C12").Value = 3.203045367113959
$ws.Range("C13").Value = 6.351244066982817
$ws.Range("C14").Value = 1.574464994382729
$ws.Range("C15").Value = 1.602253972453581
$ws.Range("C16").Value = 22.04689765473782
$ws.Range("C17").Value = 18.08916214633125
$ws.Range("C18").Value = 9.569646500925405
$ws.Range("C19").Value = 1.074994678004
$ws.Range("C20").Value = 21.06331408880899
$ws.Range("C21").Value = 63.47587622499809
$ws.Range("C22").Value = 6.501158290786097
$ws.Range("C23").Value = 0.1901351131163538
$ws.Range("C24").Value = 2.032251843655181
$ws.Range("C25").Value = 25.05395759756261
$ws.Range("C26").Value = 6.332230555671181
$ws.Range("C27").Value = 0.5433476501748109
$ws.Range("C28").Value = 9.48042925554004
$ws.Range("C29").Value = 19.50493745015149
$ws.Range("C30").Value = 8.945125783227843
$ws.Range("C31").Value = 4.84844538446702
$ws.Range("C32").Value = 3.395374346920117
$ws.Range("C33").Value = 1.782151041017516
$ws.Range("C34").Value = 4.634177737762822
$ws.Range("C35").Value = 2.302097446501083
$ws.Range("C36").Value = 80.3562178252473
$ws.Range("C37").Value = 4.75118396121904
$ws.Range("C38").Value = 25.69602924877861
$ws.Range("C39").Value = 4.68171151604191
$ws.Range("C40").Value = 3.276174256774096
$ws.Range("C41").Value = 11.89148874763472
$ws.Range("C42").Value = 0.7685846303280299
$ws.Range("C43").Value = 5.579002992171779
$ws.Range("C44").Value = 291.794968708758
